$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 216, shifting the existing rows 216-228 down to 217-229.
$ws.Rows.Item(216).Insert()

# Populate the newly inserted row 216 with the new weekly price record
# (Feria Lagunitas de Puerto Montt - Pepino ensalada).
$ws.Cells.Item(216, 1).Value = 4
$ws.Cells.Item(216, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(216, 3).Value = "Los Lagos"
$ws.Cells.Item(216, 4).Value = 44610
$ws.Cells.Item(216, 5).Value = 10
$ws.Cells.Item(216, 6).Value = 100112043
$ws.Cells.Item(216, 7).Value = "Pepino ensalada"
$ws.Cells.Item(216, 8).Value = "Sin especificar"
$ws.Cells.Item(216, 9).Value = "Primera"
$ws.Cells.Item(216, 10).Value = 400
$ws.Cells.Item(216, 11).Value = 16000
$ws.Cells.Item(216, 12).Value = 16000
$ws.Cells.Item(216, 13).Value = 16000
$ws.Cells.Item(216, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(216, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(216, 16).Value = 267
$ws.Cells.Item(216, 17).Value = 60
$ws.Cells.Item(216, 18).Value = "Hortaliza"
